$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '44.257.46'
$ws.Range("E2").Value = '  +2.78%  '

$ws.Range("D3").Value = '2.266.14'
$ws.Range("E3").Value = '  +1.82%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '320.66'
$ws.Range("E5").Value = '  -0.75%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '102.64'
$ws.Range("E6").Value = '  +3.83%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.581'
$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("E8").Value = '  +0.03%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.556'
$ws.Range("E9").Value = '  -1.33%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '37.48'
$ws.Range("E10").Value = '  +1.80%  '

$ws.Range("E11").Value = '  +1.01%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.66'
$ws.Range("E12").Value = '  +0.25%  '

$ws.Range("E13").Value = '  -0.90%  '

$ws.Range("D14").Value = '2.601.32'
$ws.Range("E14").Value = '  +1.30%  '

$ws.Range("E15").Value = '  +0.24%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.51'
$ws.Range("E16").Value = '  +0.91%  '

$ws.Range("D17").Value = '2.265.57'
$ws.Range("E17").Value = '  +1.55%  '

$ws.Range("D18").Value = '44.140.35'
$ws.Range("E18").Value = '  +2.69%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.44'
$ws.Range("E19").Value = '  -4.07%  '

$ws.Range("E20").Value = '  +2.16%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.57'
$ws.Range("E21").Value = '  +0.50%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '65.90'
$ws.Range("E22").Value = '  +1.13%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.17'
$ws.Range("E23").Value = '  -1.60%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '236.28'
$ws.Range("E24").Value = '  -0.42%  '

$ws.Range("E25").Value = '  -3.07%  '

$ws.Range("E26").Value = '  +0.07%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.33'
$ws.Range("E27").Value = '  +3.44%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '39.33'
$ws.Range("E28").Value = '  +8.41%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.19'
$ws.Range("E29").Value = '  -2.51%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.27'
$ws.Range("E30").Value = '  -1.35%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '162.92'
$ws.Range("E31").Value = '  +5.56%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '20.31'
$ws.Range("E32").Value = '  -0.27%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0855'
$ws.Range("E33").Value = '  -1.47%  '

$ws.Range("E34").Value = '  +0.45%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.115'
$ws.Range("E35").Value = '  +10.64%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.99'
$ws.Range("E36").Value = '  +4.51%  '

$ws.Range("E37").Value = '  -6.64%  '

$ws.Range("E38").Value = '  -1.57%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '16.85'
$ws.Range("E39").Value = '  +21.21%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.74'
$ws.Range("E40").Value = '  +0.76%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.25'
$ws.Range("E41").Value = '  -3.93%  '

$ws.Range("E42").Value = '  -1.08%  '

$ws.Range("E43").Value = '  +0.08%  '

$ws.Range("D44").Value = '1.787.20'
$ws.Range("E44").Value = '  +3.67%  '

$ws.Range("E45").Value = '  -1.53%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '82.76'
$ws.Range("E46").Value = '  -2.63%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '75.27'
$ws.Range("E47").Value = '  +0.55%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.24'
$ws.Range("E48").Value = '  -0.75%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '105.31'
$ws.Range("E49").Value = '  +2.32%  '

$ws.Range("E50").Value = '  +6.57%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '58.67'
$ws.Range("E51").Value = '  +1.11%  '
